$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row with the
# freshly scraped values. Each pair of cells is briefly switched to Text
# format while writing so the exact original string formatting (trailing
# zeros, thousand-dot grouping, padded percent strings, ...) survives the
# COM assignment instead of being auto-coerced to a number; the format is
# cleared right after so the cell keeps its original (default) style.

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = '27.850.11'
$ws.Range("E2").Value = '  -1.28%  '
$rng.ClearFormats()

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = '1.902.49'
$ws.Range("E3").Value = '  -0.76%  '
$rng.ClearFormats()

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.64%  '
$rng.ClearFormats()

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = '312.69'
$ws.Range("E5").Value = '  -1.47%  '
$rng.ClearFormats()

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.59%  '
$rng.ClearFormats()

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = '0.4930'
$ws.Range("E7").Value = '  +1.83%  '
$rng.ClearFormats()

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = '0.3824'
$ws.Range("E8").Value = '  -0.19%  '
$rng.ClearFormats()

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = '0.07320'
$ws.Range("E9").Value = '  -0.82%  '
$rng.ClearFormats()

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = '0.9094'
$ws.Range("E10").Value = '  -3.45%  '
$rng.ClearFormats()

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = '21.07'
$ws.Range("E11").Value = '  +0.56%  '
$rng.ClearFormats()

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = '0.07614'
$ws.Range("E12").Value = '  -2.58%  '
$rng.ClearFormats()

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = '1.898.99'
$ws.Range("E13").Value = '  -1.03%  '
$rng.ClearFormats()

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = '5.481'
$ws.Range("E14").Value = '  -0.45%  '
$rng.ClearFormats()

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = '6.649'
$ws.Range("E15").Value = '  -0.01%  '
$rng.ClearFormats()

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = '91.26'
$ws.Range("E16").Value = '  -0.13%  '
$rng.ClearFormats()

$rng = $ws.Range("D18:E18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = '0.000008746'
$ws.Range("E18").Value = '  -1.13%  '
$rng.ClearFormats()

$rng = $ws.Range("D19:E19")
$rng.NumberFormat = "@"
$ws.Range("D19").Value = '0.9993'
$ws.Range("E19").Value = '  -0.55%  '
$rng.ClearFormats()

$rng = $ws.Range("D20:E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = '27.868.75'
$ws.Range("E20").Value = '  -1.25%  '
$rng.ClearFormats()

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = '14.54'
$ws.Range("E21").Value = '  -2.28%  '
$rng.ClearFormats()

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = '5.134'
$ws.Range("E22").Value = '  -0.57%  '
$rng.ClearFormats()

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = '10.80'
$ws.Range("E23").Value = '  -1.35%  '
$rng.ClearFormats()

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = '154.41'
$ws.Range("E24").Value = '  -1.24%  '
$rng.ClearFormats()

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = '1.867'
$ws.Range("E25").Value = '  -3.05%  '
$rng.ClearFormats()

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = '2.228'
$ws.Range("E26").Value = '  +5.89%  '
$rng.ClearFormats()

$rng = $ws.Range("D27:E27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = '18.38'
$ws.Range("E27").Value = '  -1.07%  '
$rng.ClearFormats()

$rng = $ws.Range("D28:E28")
$rng.NumberFormat = "@"
$ws.Range("D28").Value = '115.11'
$ws.Range("E28").Value = '  -1.11%  '
$rng.ClearFormats()

$rng = $ws.Range("D29:E29")
$rng.NumberFormat = "@"
$ws.Range("D29").Value = '4.912'
$ws.Range("E29").Value = '  -1.34%  '
$rng.ClearFormats()

$rng = $ws.Range("D30:E30")
$rng.NumberFormat = "@"
$ws.Range("D30").Value = '0.08936'
$ws.Range("E30").Value = '  +0.26%  '
$rng.ClearFormats()

$rng = $ws.Range("D31:E31")
$rng.NumberFormat = "@"
$ws.Range("D31").Value = '3.202'
$ws.Range("E31").Value = '  -4.61%  '
$rng.ClearFormats()

$rng = $ws.Range("D32:E32")
$rng.NumberFormat = "@"
$ws.Range("D32").Value = '1.241'
$ws.Range("E32").Value = '  -0.94%  '
$rng.ClearFormats()

$rng = $ws.Range("D33:E33")
$rng.NumberFormat = "@"
$ws.Range("D33").Value = '0.7720'
$ws.Range("E33").Value = '  -0.44%  '
$rng.ClearFormats()

$rng = $ws.Range("D34:E34")
$rng.NumberFormat = "@"
$ws.Range("D34").Value = '4.643'
$ws.Range("E34").Value = '  -1.35%  '
$rng.ClearFormats()

$rng = $ws.Range("D35:E35")
$rng.NumberFormat = "@"
$ws.Range("D35").Value = '0.02063'
$ws.Range("E35").Value = '  +0.46%  '
$rng.ClearFormats()

$rng = $ws.Range("D36:E36")
$rng.NumberFormat = "@"
$ws.Range("D36").Value = '2.577'
$ws.Range("E36").Value = '  -4.10%  '
$rng.ClearFormats()

$rng = $ws.Range("D37:E37")
$rng.NumberFormat = "@"
$ws.Range("D37").Value = '1.096'
$ws.Range("E37").Value = '  -0.77%  '
$rng.ClearFormats()

$rng = $ws.Range("D38:E38")
$rng.NumberFormat = "@"
$ws.Range("D38").Value = '0.5533'
$ws.Range("E38").Value = '  -0.20%  '
$rng.ClearFormats()

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = '0.05292'
$ws.Range("E39").Value = '  -0.78%  '
$rng.ClearFormats()

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = '3.006'
$ws.Range("E40").Value = '  +0.03%  '
$rng.ClearFormats()

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = '6.988'
$ws.Range("E41").Value = '  -0.94%  '
$rng.ClearFormats()

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = '8.539'
$ws.Range("E42").Value = '  +0.65%  '
$rng.ClearFormats()

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = '0.1522'
$ws.Range("E43").Value = '  -0.72%  '
$rng.ClearFormats()

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = '111.17'
$ws.Range("E44").Value = '  +3.89%  '
$rng.ClearFormats()

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = '10.67'
$ws.Range("E45").Value = '  -1.03%  '
$rng.ClearFormats()

$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = '0.4798'
$ws.Range("E46").Value = '  -1.46%  '
$rng.ClearFormats()

$rng = $ws.Range("D47:E47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = '0.9994'
$ws.Range("E47").Value = '  -0.62%  '
$rng.ClearFormats()

$rng = $ws.Range("D48:E48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = '1.642'
$ws.Range("E48").Value = '  -1.14%  '
$rng.ClearFormats()

$rng = $ws.Range("D49:E49")
$rng.NumberFormat = "@"
$ws.Range("D49").Value = '67.50'
$ws.Range("E49").Value = '  -1.68%  '
$rng.ClearFormats()

$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = '0.06070'
$ws.Range("E50").Value = '  -0.85%  '
$rng.ClearFormats()

$rng = $ws.Range("D51:E51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = '0.8979'
$ws.Range("E51").Value = '  -1.30%  '
$rng.ClearFormats()
